# Add two new "2-option trial" sheets (TAG7, TAG8) with their raw trial
# data, and append their summary rows to the "summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add TAG7 sheet (after TAG6, i.e. at the end of the tab order)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tag7 = $wb.Worksheets.Add($null, $lastSheet)
$tag7.Name = "TAG7"

$tag7.Range("A1").Value = "Trial"
$tag7.Range("B1").Value = "Reaction Time"
$tag7.Range("C1").Value = "Accuracy"

$tag7.Range("A2").Value = "('folder_4', 'ball_4')"
$tag7.Range("B2").Value = "2.51s"
$tag7.Range("C2").Value = "correct"

$tag7.Range("A3").Value = "('atest', 'ball_1 - Copy')"
$tag7.Range("B3").Value = "2.03s"
$tag7.Range("C3").Value = "correct"

$tag7.Range("A4").Value = "('folder_2', 'ball_2')"
$tag7.Range("B4").Value = "1.98s"
$tag7.Range("C4").Value = "correct"

$tag7.Range("A5").Value = "('atest', 'ball_1')"
$tag7.Range("B5").Value = "1.91s"
$tag7.Range("C5").Value = "correct"

$tag7.Range("A6").Value = "('folder_3', 'ball_3')"
$tag7.Range("B6").Value = "1.88s"
$tag7.Range("C6").Value = "correct"

# ---------------------------------------------------------------------
# 2. Add TAG8 sheet (after TAG7, i.e. at the end of the tab order)
# ---------------------------------------------------------------------
$tag8 = $wb.Worksheets.Add($null, $tag7)
$tag8.Name = "TAG8"

$tag8.Range("A1").Value = "Trial"
$tag8.Range("B1").Value = "Reaction Time"
$tag8.Range("C1").Value = "Accuracy"

$tag8.Range("A2").Value = "('atest', 'ball_1 - Copy')"
$tag8.Range("B2").Value = "2.81s"
$tag8.Range("C2").Value = "correct"

$tag8.Range("A3").Value = "('folder_3', 'ball_3')"
$tag8.Range("B3").Value = "2.11s"
$tag8.Range("C3").Value = "correct"

$tag8.Range("A4").Value = "('atest', 'ball_1')"
$tag8.Range("B4").Value = "2.03s"
$tag8.Range("C4").Value = "correct"

$tag8.Range("A5").Value = "('folder_2', 'ball_2')"
$tag8.Range("B5").Value = "1.08s"
$tag8.Range("C5").Value = "correct"

$tag8.Range("A6").Value = "('folder_4', 'ball_4')"
$tag8.Range("B6").Value = "1.87s"
$tag8.Range("C6").Value = "correct"

# ---------------------------------------------------------------------
# 3. Append the summary rows for TAG7 / TAG8 onto the "summary" sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("summary")

$summary.Range("A8").Value = "TAG7"
$summary.Range("B8").Value = "2.06s"
$summary.Range("C8").Value = "5/5"
$summary.Range("D8").Value = "2.06s"
$summary.Range("E8").Value = "5/5"

$summary.Range("A9").Value = "TAG8"
$summary.Range("B9").Value = "1.98s"
$summary.Range("C9").Value = "5/5"
$summary.Range("D9").Value = "1.98s"
$summary.Range("E9").Value = "5/5"
